$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("extraction")

# Clear the "imageTemplate" values for the two data rows (E2, E3).
# The field was present but empty -> clear contents (keeps cell, removes
# hyperlink/text) rather than deleting the cell entirely.
$ws.Range("E2:E3").ClearContents()
$ws.Range("E2:E3").Hyperlinks.Delete()

$ws.Range("E3").Select()
